# Update the A/D column values (shared strings TestSignupb# -> TestSignupc#)
# and move the sheet selection to D8, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TestSignupc9"
$ws.Range("A3").Value = "TestSignupc10"
$ws.Range("A4").Value = "TestSignupc11"
$ws.Range("A5").Value = "TestSignupc12"

$ws.Range("D2").Value = "TestSignupc9"
$ws.Range("D3").Value = "TestSignupc10"
$ws.Range("D4").Value = "TestSignupc11"
$ws.Range("D5").Value = "TestSignupc12"

$ws.Range("D8").Select()
